$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "보고서 날짜: 2024년 1월 22일" "보고 날짜: 2024년 1월 22일"

Replace-Text " 바이럴 릴은 매력적인 시각적 개체와 유익한 콘텐츠를 통해 제품의 매력을 성공적으로 전달했습니다. " " 바이럴 쇼츠는 매력적인 시각적 개체와 유익한 콘텐츠를 통해 제품의 매력을 성공적으로 전달했습니다. "

Replace-Text " 인플루언서 마케팅의 힘은 과장될 수 없습니다. " " 인플루언서 마케팅의 힘은 아무리 강조해도 지나치지 않습니다. "

Replace-Text "맛과 맛 품종:" "맛과 풍미의 다양성:"

Replace-Text " 맛있고 다양한 맛에 대한 Contoso Protein Plus의 명성은 바이러스 성 콘텐츠의 주요 판매 포인트였습니다. " " 맛과 다양한 풍미에 대한 Contoso Protein Plus의 명성은 바이럴 콘텐츠의 주요 판매 포인트였습니다. "

Replace-Text "건강 및 피트니스 트렌드:" "건강 및 피트니스 추세:"

Replace-Text " 건강과 피트니스 의식의 지속적인 급증, 운동 루틴과 활동적인 라이프 스타일을 포용하는 사람들의 수의 증가와 함께, Contoso 단백질 플러스 같은 제품에 대한 수용 시장을 만들었습니다." " 건강과 피트니스 의식의 지속적인 급증, 운동 루틴과 활동적인 라이프 스타일을 포용하는 사람들의 수적 증가와 함께, Contoso Protein Plus 같은 제품을 수용하는 시장이 형성되었습니다."

Replace-Text " 다양한 온라인 소매 업체를 통한 제품의 접근성은 과대 광고를 더욱 촉진시켰습니다. " " 다양한 온라인 소매업체를 통한 제품의 접근성으로 인해 광고가 더욱 촉진되었습니다. "

Replace-Text "긍정적 인 리뷰 및 평가 :" "긍정적인 리뷰 및 평가:"

Replace-Text " 릴은 고립 된 경우가 아니었습니다. " " 쇼츠는 격리된 사례가 아니었습니다. "
